$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -20.99795418726239
$ws.Range("F2").Value = -45.59603536128998
$ws.Range("B3").Value = -21.10179198563014
$ws.Range("F3").Value = -41.06186187267303
$ws.Range("B4").Value = -21.55464948246708
$ws.Range("F4").Value = -37.75727462768555
$ws.Range("B5").Value = -21.92792841161213
$ws.Range("F5").Value = -36.11608040332794
$ws.Range("B6").Value = -22.11964205490517
$ws.Range("F6").Value = -35.42429411411285
$ws.Range("B7").Value = -22.11964205490517
$ws.Range("F7").Value = -35.42429411411285
$ws.Range("B8").Value = -22.07057130856833
$ws.Range("F8").Value = -35.59410583972931
$ws.Range("B9").Value = -21.48242981463659
$ws.Range("F9").Value = -38.14239275455475
$ws.Range("B10").Value = -20.98795145568283
$ws.Range("F10").Value = -45.24657070636749
$ws.Range("B11").Value = -21.79815971245625
$ws.Range("F11").Value = -53.13354885578156
$ws.Range("B12").Value = -22.62771894779121
$ws.Range("F12").Value = -57.27001345157623
$ws.Range("B13").Value = -23.47091593793789
$ws.Range("F13").Value = -60.65550637245178
$ws.Range("B14").Value = -23.6366354565389
$ws.Range("F14").Value = -61.26440155506134
$ws.Range("B15").Value = -23.39948012117884
$ws.Range("F15").Value = -60.38825237751007
$ws.Range("B16").Value = -23.40169217802895
$ws.Range("F16").Value = -60.39657354354858
$ws.Range("B17").Value = -23.41583489571281
$ws.Range("F17").Value = -60.44970417022705
$ws.Range("B18").Value = -22.72146045022919
$ws.Range("F18").Value = -57.67579197883606
$ws.Range("B19").Value = -22.35007561258658
$ws.Range("F19").Value = -56.00964665412903
$ws.Range("B20").Value = -22.35265164940233
$ws.Range("F20").Value = -56.02179384231567
$ws.Range("B21").Value = -23.41720580255719
$ws.Range("F21").Value = -60.45484805107117
$ws.Range("B22").Value = -24.74878119468531
$ws.Range("F22").Value = -65.03224408626556
$ws.Range("B23").Value = -25.97044579442127
$ws.Range("F23").Value = -68.72418415546417
$ws.Range("B24").Value = -26.47008051876992
$ws.Range("F24").Value = -70.13933300971985
$ws.Range("B25").Value = -24.12734139211216
$ws.Range("F25").Value = -62.98786759376526
